$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Cells.Item(17, 8).Value = 675407.5  # H17
$ws.Cells.Item(17, 10).Value = 675407.5  # J17
$ws.Cells.Item(17, 12).Value = 2026222.5  # L17
$ws.Cells.Item(17, 14).Value = -2026558.5  # N17
# Row 64
$ws.Cells.Item(64, 8).Value = 7649.75  # H64
$ws.Cells.Item(64, 10).Value = 12374.625  # J64
$ws.Cells.Item(64, 12).Value = 12374.625  # L64
$ws.Cells.Item(64, 14).Value = -12870.625  # N64
# Row 67
$ws.Cells.Item(67, 8).Value = 7649.75  # H67
$ws.Cells.Item(67, 10).Value = 12374.625  # J67
$ws.Cells.Item(67, 12).Value = 12374.625  # L67
$ws.Cells.Item(67, 14).Value = -14090.625  # N67
# Row 86
$ws.Cells.Item(86, 8).Value = 5415.8335  # H86
$ws.Cells.Item(86, 9).Value = 4750  # I86
$ws.Cells.Item(86, 10).Value = 5748.75  # J86
$ws.Cells.Item(86, 11).Value = 4750  # K86
$ws.Cells.Item(86, 12).Value = 5748.75  # L86
$ws.Cells.Item(86, 13).Value = -3627  # M86
$ws.Cells.Item(86, 14).Value = -7994.75  # N86
# Row 89
$ws.Cells.Item(89, 8).Value = 5415.8335  # H89
$ws.Cells.Item(89, 9).Value = 4750  # I89
$ws.Cells.Item(89, 10).Value = 5748.75  # J89
$ws.Cells.Item(89, 11).Value = 23750  # K89
$ws.Cells.Item(89, 12).Value = 28743.75  # L89
$ws.Cells.Item(89, 13).Value = -18134  # M89
$ws.Cells.Item(89, 14).Value = -39975.75  # N89
# Row 107
$ws.Cells.Item(107, 8).Value = 788.17645  # H107
$ws.Cells.Item(107, 9).Value = 768.1429000000001  # I107
$ws.Cells.Item(107, 10).Value = 881.6667  # J107
$ws.Cells.Item(107, 11).Value = 768.1429000000001  # K107
$ws.Cells.Item(107, 12).Value = 881.6667  # L107
$ws.Cells.Item(107, 13).Value = 1151.8571  # M107
$ws.Cells.Item(107, 14).Value = -4721.6667  # N107
# Row 129
$ws.Cells.Item(129, 8).Value = 200000800  # H129
$ws.Cells.Item(129, 9).Value = 500000500  # I129
$ws.Cells.Item(129, 10).Value = 1000  # J129
$ws.Cells.Item(129, 11).Value = 1500001500  # K129
$ws.Cells.Item(129, 12).Value = 3000  # L129
$ws.Cells.Item(129, 13).Value = -1499996500  # M129
$ws.Cells.Item(129, 14).Value = -13000  # N129
# Row 138
$ws.Cells.Item(138, 8).Value = 5729.69  # H138
$ws.Cells.Item(138, 10).Value = 5995.7046  # J138
$ws.Cells.Item(138, 12).Value = 17987.1138  # L138
$ws.Cells.Item(138, 14).Value = -28267.1138  # N138

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 11701.3  # H32
$ws.Cells.Item(32, 9).Value = 10522.448  # I32
$ws.Cells.Item(32, 10).Value = 39993.75  # J32
$ws.Cells.Item(32, 11).Value = 10522.448  # K32
$ws.Cells.Item(32, 12).Value = 39993.75  # L32
$ws.Cells.Item(32, 13).Value = -10235.448  # M32
$ws.Cells.Item(32, 14).Value = -40567.75  # N32
# Row 61
$ws.Cells.Item(61, 8).Value = 538365.75  # H61
$ws.Cells.Item(61, 9).Value = 3941.9092  # I61
$ws.Cells.Item(61, 11).Value = 3941.9092  # K61
$ws.Cells.Item(61, 13).Value = -3729.9092  # M61
# Row 63
$ws.Cells.Item(63, 8).Value = 4911.1665  # H63
$ws.Cells.Item(63, 9).Value = 4866.75  # I63
$ws.Cells.Item(63, 11).Value = 4866.75  # K63
$ws.Cells.Item(63, 13).Value = -4180.75  # M63
# Row 66
$ws.Cells.Item(66, 8).Value = 4911.1665  # H66
$ws.Cells.Item(66, 9).Value = 4866.75  # I66
$ws.Cells.Item(66, 11).Value = 24333.75  # K66
$ws.Cells.Item(66, 13).Value = -20901.75  # M66
# Row 74
$ws.Cells.Item(74, 8).Value = 48374.76  # H74
$ws.Cells.Item(74, 9).Value = 51343.8  # I74
$ws.Cells.Item(74, 10).Value = 36498.6  # J74
$ws.Cells.Item(74, 11).Value = 51343.8  # K74
$ws.Cells.Item(74, 12).Value = 36498.6  # L74
$ws.Cells.Item(74, 13).Value = -50469.8  # M74
$ws.Cells.Item(74, 14).Value = -38246.6  # N74
# Row 77
$ws.Cells.Item(77, 8).Value = 48374.76  # H77
$ws.Cells.Item(77, 9).Value = 51343.8  # I77
$ws.Cells.Item(77, 10).Value = 36498.6  # J77
$ws.Cells.Item(77, 11).Value = 256719  # K77
$ws.Cells.Item(77, 12).Value = 182493  # L77
$ws.Cells.Item(77, 13).Value = -252351  # M77
$ws.Cells.Item(77, 14).Value = -191229  # N77
# Row 110
$ws.Cells.Item(110, 8).Value = 9599.538  # H110
$ws.Cells.Item(110, 9).Value = 11873.883  # I110
$ws.Cells.Item(110, 11).Value = 11873.883  # K110
$ws.Cells.Item(110, 13).Value = -9828.883  # M110
# Row 122
$ws.Cells.Item(122, 8).Value = 2712  # H122
$ws.Cells.Item(122, 9).Value = 2712  # I122
$ws.Cells.Item(122, 11).Value = 8136  # K122
$ws.Cells.Item(122, 13).Value = -5686  # M122
# Row 132
$ws.Cells.Item(132, 8).Value = 5830.222  # H132
$ws.Cells.Item(132, 9).Value = 5429.5884  # I132
$ws.Cells.Item(132, 11).Value = 16288.7652  # K132
$ws.Cells.Item(132, 13).Value = -13758.7652  # M132
# Row 136
$ws.Cells.Item(136, 8).Value = 538365.75  # H136
$ws.Cells.Item(136, 9).Value = 3941.9092  # I136
$ws.Cells.Item(136, 11).Value = 11825.7276  # K136
$ws.Cells.Item(136, 13).Value = -9275.7276  # M136
# Row 139
$ws.Cells.Item(139, 8).Value = 100178.75  # H139
$ws.Cells.Item(139, 10).Value = 100178.75  # J139
$ws.Cells.Item(139, 12).Value = 100178.75  # L139
$ws.Cells.Item(139, 14).Value = -110458.75  # N139

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Cells.Item(105, 8).Value = 1242.6666  # H105
$ws.Cells.Item(105, 9).Value = 1242.6666  # I105
$ws.Cells.Item(105, 11).Value = 1242.6666  # K105
$ws.Cells.Item(105, 13).Value = 504.3334  # M105

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 2801.0715  # H31
$ws.Cells.Item(31, 9).Value = 2537.36  # I31
$ws.Cells.Item(31, 11).Value = 2537.36  # K31
$ws.Cells.Item(31, 13).Value = -2242.36  # M31
# Row 34
$ws.Cells.Item(34, 8).Value = 2801.0715  # H34
$ws.Cells.Item(34, 9).Value = 2537.36  # I34
$ws.Cells.Item(34, 11).Value = 2537.36  # K34
$ws.Cells.Item(34, 13).Value = -2335.36  # M34
# Row 58
$ws.Cells.Item(58, 8).Value = 3426  # H58
$ws.Cells.Item(58, 9).Value = 3799.8  # I58
$ws.Cells.Item(58, 11).Value = 3799.8  # K58
$ws.Cells.Item(58, 13).Value = -3596.8  # M58
# Row 99
$ws.Cells.Item(99, 8).Value = 5185.875  # H99
$ws.Cells.Item(99, 10).Value = 6250  # J99
$ws.Cells.Item(99, 12).Value = 6250  # L99
$ws.Cells.Item(99, 14).Value = -9246  # N99
# Row 105
$ws.Cells.Item(105, 8).Value = 990  # H105
$ws.Cells.Item(105, 9).Value = 990  # I105
$ws.Cells.Item(105, 10).Value = 0  # J105
$ws.Cells.Item(105, 11).Value = 990  # K105
$ws.Cells.Item(105, 12).Value = 0  # L105
$ws.Cells.Item(105, 13).ClearContents()  # M105
$ws.Cells.Item(105, 14).Value = 757  # N105
# Row 126
$ws.Cells.Item(126, 8).Value = 5185.875  # H126
$ws.Cells.Item(126, 10).Value = 6250  # J126
$ws.Cells.Item(126, 12).Value = 18750  # L126
$ws.Cells.Item(126, 14).Value = -23690  # N126
# Row 132
$ws.Cells.Item(132, 8).Value = 3493.7856  # H132
$ws.Cells.Item(132, 9).Value = 3181.7273  # I132
$ws.Cells.Item(132, 10).Value = 4638  # J132
$ws.Cells.Item(132, 11).Value = 9545.1819  # K132
$ws.Cells.Item(132, 12).Value = 13914  # L132
$ws.Cells.Item(132, 13).Value = -7015.1819  # M132
$ws.Cells.Item(132, 14).Value = -18974  # N132
# Row 134
$ws.Cells.Item(134, 8).Value = 20232.158  # H134
$ws.Cells.Item(134, 9).Value = 14993.643  # I134
$ws.Cells.Item(134, 11).Value = 44980.929  # K134
$ws.Cells.Item(134, 13).Value = -42445.929  # M134
# Row 136
$ws.Cells.Item(136, 8).Value = 3426  # H136
$ws.Cells.Item(136, 9).Value = 3799.8  # I136
$ws.Cells.Item(136, 11).Value = 11399.4  # K136
$ws.Cells.Item(136, 13).Value = -8849.400000000001  # M136

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 22
$ws.Cells.Item(22, 8).Value = 104666.664  # H22
$ws.Cells.Item(22, 9).Value = 300000  # I22
$ws.Cells.Item(22, 10).Value = 7000  # J22
$ws.Cells.Item(22, 11).Value = 900000  # K22
$ws.Cells.Item(22, 12).Value = 21000  # L22
$ws.Cells.Item(22, 13).Value = -899831  # M22
$ws.Cells.Item(22, 14).Value = -21338  # N22
# Row 27
$ws.Cells.Item(27, 8).Value = 104666.664  # H27
$ws.Cells.Item(27, 9).Value = 300000  # I27
$ws.Cells.Item(27, 10).Value = 7000  # J27
$ws.Cells.Item(27, 11).Value = 900000  # K27
$ws.Cells.Item(27, 12).Value = 21000  # L27
$ws.Cells.Item(27, 13).Value = -899898  # M27
$ws.Cells.Item(27, 14).Value = -21204  # N27
# Row 137
$ws.Cells.Item(137, 8).Value = 4130.7144  # H137
$ws.Cells.Item(137, 10).Value = 4500  # J137
$ws.Cells.Item(137, 12).Value = 13500  # L137
$ws.Cells.Item(137, 14).Value = -23700  # N137

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 46
$ws.Cells.Item(46, 8).Value = 50000  # H46
$ws.Cells.Item(46, 10).Value = 50000  # J46
$ws.Cells.Item(46, 12).Value = 50000  # L46
$ws.Cells.Item(46, 14).Value = -50312  # N46
# Row 102
$ws.Cells.Item(102, 8).Value = 40001436  # H102
$ws.Cells.Item(102, 9).Value = 1564.3  # I102
$ws.Cells.Item(102, 11).Value = 1564.3  # K102
$ws.Cells.Item(102, 13).Value = 57.70000000000005  # M102
# Row 132
$ws.Cells.Item(132, 8).Value = 4845.6  # H132
$ws.Cells.Item(132, 9).Value = 3200.1667  # I132
$ws.Cells.Item(132, 10).Value = 5942.5557  # J132
$ws.Cells.Item(132, 11).Value = 9600.500100000001  # K132
$ws.Cells.Item(132, 12).Value = 17827.6671  # L132
$ws.Cells.Item(132, 13).Value = -7070.500100000001  # M132
$ws.Cells.Item(132, 14).Value = -22887.6671  # N132
# Row 136
$ws.Cells.Item(136, 8).Value = 0  # H136
$ws.Cells.Item(136, 10).Value = 0  # J136
$ws.Cells.Item(136, 12).ClearContents()  # L136
$ws.Cells.Item(136, 14).Value = 0  # N136
# Row 137
$ws.Cells.Item(137, 8).Value = 70000  # H137
$ws.Cells.Item(137, 10).Value = 0  # J137
$ws.Cells.Item(137, 12).Value = 0  # L137
$ws.Cells.Item(137, 14).ClearContents()  # N137

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Cells.Item(40, 8).Value = 8489.299999999999  # H40
$ws.Cells.Item(40, 9).Value = 8320  # I40
$ws.Cells.Item(40, 11).Value = 8320  # K40
$ws.Cells.Item(40, 13).Value = -8184  # M40
# Row 122
$ws.Cells.Item(122, 8).Value = 424924.6  # H122
$ws.Cells.Item(122, 9).Value = 674446.4399999999  # I122
$ws.Cells.Item(122, 11).Value = 2023339.32  # K122
$ws.Cells.Item(122, 13).Value = -2020889.32  # M122
# Row 136
$ws.Cells.Item(136, 8).Value = 5132.0835  # H136
$ws.Cells.Item(136, 9).Value = 4758.5  # I136
$ws.Cells.Item(136, 10).Value = 7000  # J136
$ws.Cells.Item(136, 11).Value = 14275.5  # K136
$ws.Cells.Item(136, 12).Value = 21000  # L136
$ws.Cells.Item(136, 13).Value = -11725.5  # M136
$ws.Cells.Item(136, 14).Value = -26100  # N136

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Cells.Item(132, 8).Value = 4347.923  # H132
$ws.Cells.Item(132, 9).Value = 4390.7646  # I132
$ws.Cells.Item(132, 13).Value = -10642.2938  # M132
# Row 137
$ws.Cells.Item(137, 8).Value = 107166.664  # H137
$ws.Cells.Item(137, 10).Value = 107166.664  # J137
$ws.Cells.Item(137, 12).Value = 107166.664  # L137
$ws.Cells.Item(137, 14).Value = -117366.664  # N137
